# Apply the two notes-page text edits described by the commit.
#
# Slide 4 notes ("Notes Placeholder 2"): tighten the presenter-timing
# sentence.
#
# Slide 7 notes ("Notes Placeholder 2"): merge the trailing "page changes"
# run back into the main sentence so the whole note is one contiguous
# run/paragraph.

$p = $ppt.ActivePresentation

# --- Slide 4 notes: update "Presenter Nikita ad Veera ..." paragraph ---
$slide4 = $p.Slides.Item(4)
$notes4 = $slide4.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes4.Paragraphs(1).Text = "Presenter Nikita ad Veera 6 minutes max for 4, 5 and 6."

# --- Slide 7 notes: merge the final paragraph's two runs into one ---
$slide7 = $p.Slides.Item(7)
$notes7 = $slide7.NotesPage.Shapes.Item(2).TextFrame.TextRange
$lastParaIndex = $notes7.Paragraphs().Count
$notes7.Paragraphs($lastParaIndex).Text = "At the end of presentation. Jamie brings CLI and GUI on sharing screen >>  Kalyan makes change to insert value file in GIT >> Shows the pods rebuilding and show web page changes"
